function Set-ParagraphXml {
    param($doc, $paraIndex, $innerXml)
    $p = $doc.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $delRange = $doc.Range($r.Start, $r.End - 1)
    $delRange.Text = ""
    $p2 = $doc.Paragraphs.Item($paraIndex)
    $r2 = $p2.Range
    $collapsed = $doc.Range($r2.Start, $r2.Start)
    $wrapper = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $collapsed.InsertXML($wrapper)
}

$d = $word.ActiveDocument

# --- Paragraph: "Backup panel\'s background color the same as ..." ---
$para1Xml = '<w:r><w:t xml:space="preserve"> In the </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t>Backup and Restore</w:t></w:r><w:r><w:t xml:space="preserve"> screen, make the </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t>Backup panel’s background</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t xml:space="preserve"> (Panel where the file name is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t>shon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t>color</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the same as the </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t xml:space="preserve">Available Backups panel’s background </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t>color</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r>'
Set-ParagraphXml $d 208 $para1Xml

# --- Paragraph: "Additionally, set a white border for each backup panel." ---
$para2Xml = '<w:r><w:t xml:space="preserve">Additionally, set a </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Strong"/><w:rFonts w:eastAsiaTheme="majorEastAsia"/></w:rPr><w:t>white border</w:t></w:r><w:r><w:t xml:space="preserve"> for each backup panel</w:t></w:r><w:r><w:t xml:space="preserve"> where file name is shown</w:t></w:r><w:r><w:t>.</w:t></w:r>'
Set-ParagraphXml $d 209 $para2Xml

# --- Paragraph: "The Splash Screen, having text ..." ---
$para3Xml = '<w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Splash Screen</w:t></w:r><w:r><w:t xml:space="preserve"> displaying the text </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>“My Records – Secure, Organised, Accessible”</w:t></w:r><w:r><w:t xml:space="preserve"> briefly appears (for about one second) when interacting with certain UI elements such as the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>search text box</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>dark/light mode toggle</w:t></w:r><w:r><w:t xml:space="preserve">, or during the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Change PIN</w:t></w:r><w:r><w:t xml:space="preserve"> operation.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">This issue should be </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>fixed</w:t></w:r><w:r><w:t xml:space="preserve"> so that the splash screen only appears </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>once at app launch</w:t></w:r><w:r><w:t xml:space="preserve"> and not during normal app interactions.</w:t></w:r>'
Set-ParagraphXml $d 212 $para3Xml

# --- Paragraph: "Build and install the debug APK ..." -> add magenta highlight ---
$p4 = $d.Paragraphs.Item(214)
$p4.Range.HighlightColorIndex = 5

Write-Host "Done"
